$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.044.70'
$ws.Range('E2').Value = '  -2.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.127.52'
$ws.Range('E3').Value = '  -0.79%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.56'
$ws.Range('E5').Value = '  -2.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.10'
$ws.Range('E6').Value = '  -5.79%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.121.37'
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  -1.78%  '
$ws.Range('E10').Value = '  -3.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.29'
$ws.Range('E11').Value = '  -1.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.460'
$ws.Range('E12').Value = '  -2.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000249'
$ws.Range('E13').Value = '  -3.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.04'
$ws.Range('E14').Value = '  -4.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.640.39'
$ws.Range('E15').Value = '  -0.79%  '
$ws.Range('E16').Value = '  +0.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.037.38'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.127.27'
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.71'
$ws.Range('E19').Value = '  -2.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '475.55'
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.16'
$ws.Range('E21').Value = '  -3.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.698'
$ws.Range('E22').Value = '  -4.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.64'
$ws.Range('E23').Value = '  -3.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.44'
$ws.Range('E24').Value = '  +3.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.05'
$ws.Range('E25').Value = '  -5.19%  '
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('E27').Value = '  -3.53%  '
$ws.Range('E28').Value = '  -4.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.02'
$ws.Range('E29').Value = '  -7.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.07'
$ws.Range('E30').Value = '  -1.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '27.24'
$ws.Range('E31').Value = '  +1.92%  '
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.108'
$ws.Range('E33').Value = '  -10.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.54'
$ws.Range('E34').Value = '  -4.61%  '
$ws.Range('E35').Value = '  -3.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.85'
$ws.Range('E36').Value = '  -1.98%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.98'
$ws.Range('E37').Value = '  -1.53%  '
$ws.Range('E38').Value = '  -4.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0391'
$ws.Range('E39').Value = '  -1.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '422.10'
$ws.Range('E40').Value = '  -7.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.117'
$ws.Range('E41').Value = '  -0.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.27'
$ws.Range('E42').Value = '  -1.09%  '
$ws.Range('E43').Value = '  -12.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.863.33'
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.258'
$ws.Range('E45').Value = '  -4.75%  '
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.11'
$ws.Range('E47').Value = '  -7.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.114'
$ws.Range('E48').Value = '  -0.61%  '
$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.29'
$ws.Range('E49').Value = '  -6.42%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.39'
$ws.Range('E50').Value = '  -4.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '118.16'
$ws.Range('E51').Value = '  -1.85%  '
